$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.289.21'
$ws.Range('E2').Value = '  +1.74%  '

$ws.Range('D3').Value = '2.097.91'
$ws.Range('E3').Value = '  +4.72%  '

$ws.Range('E4').Value = '  +0.23%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '250.37'
$ws.Range('E5').Value = '  +1.77%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.655'
$ws.Range('E6').Value = '  -0.78%  '

$ws.Range('E7').Value = '  +0.08%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '51.27'
$ws.Range('E8').Value = '  +13.23%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '61.63'
$ws.Range('E9').Value = '  +2.00%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.371'
$ws.Range('E10').Value = '  +2.89%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0741'
$ws.Range('E11').Value = '  +3.09%  '

$ws.Range('E12').Value = '  +6.09%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '15.31'
$ws.Range('E13').Value = '  +5.30%  '

$ws.Range('D14').Value = '2.400.06'
$ws.Range('E14').Value = '  +3.50%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.828'
$ws.Range('E15').Value = '  +3.16%  '

$ws.Range('D16').Value = '2.095.19'
$ws.Range('E16').Value = '  +4.61%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.08'
$ws.Range('E17').Value = '  +3.69%  '

$ws.Range('D18').Value = '37.264.01'
$ws.Range('E18').Value = '  +2.59%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '72.08'
$ws.Range('E19').Value = '  +1.21%  '

$ws.Range('D20').Value = '0.0₃0826'
$ws.Range('E20').Value = '  +0.78%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.49'
$ws.Range('E21').Value = '  +4.56%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '239.86'
$ws.Range('E22').Value = '  +1.93%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.21'
$ws.Range('E23').Value = '  +4.46%  '

$ws.Range('E24').Value = '  +0.24%  '

$ws.Range('E25').Value = '  +0.18%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '169.57'
$ws.Range('E26').Value = '  +4.66%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.15'
$ws.Range('E27').Value = '  +7.93%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.55'
$ws.Range('E28').Value = '  +3.31%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.99'
$ws.Range('E29').Value = '  +1.08%  '

$ws.Range('E30').Value = '  +0.24%  '

$ws.Range('E31').Value = '  +24.38%  '

$ws.Range('E32').Value = '  +3.07%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0606'
$ws.Range('E33').Value = '  +4.43%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0923'
$ws.Range('E34').Value = '  +14.55%  '

$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.39%  '

$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.32'
$ws.Range('E36').Value = '  +10.99%  '

$ws.Range('B37').Value = 'Gas'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '19.17'
$ws.Range('E37').Value = '  -10.00%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.07'
$ws.Range('E38').Value = '  +0.61%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.82'
$ws.Range('E39').Value = '  -2.24%  '

$ws.Range('E40').Value = '  -1.76%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '17.83'
$ws.Range('E41').Value = '  +10.62%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0223'
$ws.Range('E42').Value = '  +2.62%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.16'
$ws.Range('E43').Value = '  +9.18%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '98.58'
$ws.Range('E44').Value = '  +2.05%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.73'
$ws.Range('E45').Value = '  -2.99%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0894'
$ws.Range('E46').Value = '  +9.18%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.03'
$ws.Range('E47').Value = '  +10.01%  '

$ws.Range('D48').Value = '1.318.92'
$ws.Range('E48').Value = '  +1.09%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.95'
$ws.Range('E49').Value = '  +13.82%  '

$ws.Range('D50').Value = '2.296.64'
$ws.Range('E50').Value = '  +3.81%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.28'
$ws.Range('E51').Value = '  +2.54%  '
